$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.616.77"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "2.272.21"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.35"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.49"
$ws.Range("E6").Value = "  +2.15%  "

$ws.Range("E7").Value = "  -0.65%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.77"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("E11").Value = "  +2.56%  "

$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.79"
$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").Value = "2.624.27"
$ws.Range("E14").Value = "  +0.93%  "

$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("D16").Value = "2.270.11"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.790"
$ws.Range("E17").Value = "  -0.89%  "

$ws.Range("D18").Value = "42.524.99"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.32"
$ws.Range("E19").Value = "  -2.11%  "

$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.46"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.12"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.51"
$ws.Range("E27").Value = "  +1.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.59"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.58"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.43"
$ws.Range("E31").Value = "  +4.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.24"
$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.73"
$ws.Range("E34").Value = "  +2.84%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.08"
$ws.Range("E35").Value = "  -3.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0728"
$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  -1.49%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("E41").Value = "  +2.14%  "

$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.26"
$ws.Range("E42").Value = "  -7.37%  "

$ws.Range("D43").Value = "1.946.33"
$ws.Range("E43").Value = "  -2.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.88"
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.75"
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.34"
$ws.Range("E48").Value = "  +1.73%  "

$ws.Range("D49").Value = "2.494.50"
$ws.Range("E49").Value = "  +0.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.90"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.98"
$ws.Range("E51").Value = "  +0.43%  "
